$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.439.90"
$ws.Range("E2").Value = "  +3.76%  "
$ws.Range("D3").Value = "3.193.87"
$ws.Range("E3").Value = "  +5.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'205.89"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "'636.32"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D8").Value = "'0.236"
$ws.Range("E8").Value = "  +12.45%  "
$ws.Range("E9").Value = "  +5.71%  "
$ws.Range("D10").Value = "3.193.34"
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  +33.12%  "
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E13").Value = "  +8.83%  "
$ws.Range("D14").Value = "3.779.32"
$ws.Range("D15").Value = "'0.0000228"
$ws.Range("E15").Value = "  +16.83%  "
$ws.Range("D16").Value = "'31.79"
$ws.Range("E16").Value = "  +7.69%  "
$ws.Range("D17").Value = "79.236.35"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "3.194.09"
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("E19").Value = "  +8.27%  "
$ws.Range("D20").Value = "'3.11"
$ws.Range("E20").Value = "  +34.67%  "
$ws.Range("D21").Value = "'9.16"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "'429.30"
$ws.Range("E22").Value = "  +14.01%  "
$ws.Range("D23").Value = "'5.02"
$ws.Range("E23").Value = "  +14.90%  "
$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "  +12.57%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'4.80"
$ws.Range("E25").Value = "  +9.40%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.359.43"
$ws.Range("E26").Value = "  +5.51%  "
$ws.Range("D27").Value = "'76.76"
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'0.0000119"
$ws.Range("E29").Value = "  +6.08%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'9.02"
$ws.Range("E30").Value = "  +8.11%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("D33").Value = "'529.14"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("E35").Value = "  +27.37%  "
$ws.Range("D36").Value = "'22.97"
$ws.Range("E36").Value = "  +10.09%  "
$ws.Range("D37").Value = "'0.121"
$ws.Range("E37").Value = "  +12.50%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").Value = "'164.73"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'194.16"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'20.02"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D44").Value = "'5.48"
$ws.Range("E44").Value = "  +4.81%  "
$ws.Range("E45").Value = "  +11.05%  "
$ws.Range("D46").Value = "'1.80"
$ws.Range("E46").Value = "  +7.84%  "
$ws.Range("D47").Value = "'1.32"
$ws.Range("E47").Value = "  +4.65%  "
$ws.Range("D48").Value = "'43.17"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "'26.07"
$ws.Range("E49").Value = "  +15.20%  "
$ws.Range("D50").Value = "'2.55"
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").Value = "'0.634"
$ws.Range("E51").Value = "  +4.33%  "
